$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 95
$ws.Range("L8").Value = 181
$ws.Range("L9").Value = 18
$ws.Range("L10").Value = 15
$ws.Range("L15").Value = 21
$ws.Range("L19").Value = 94
$ws.Range("L20").Value = 76
$ws.Range("L21").Value = 8
$ws.Range("L22").Value = 8
$ws.Range("L24").Value = 5
$ws.Range("L27").Value = 31
$ws.Range("L29").Value = 134
$ws.Range("L33").Value = 125
$ws.Range("L37").Value = 95
$ws.Range("L42").Value = 93
$ws.Range("L48").Value = 49
$ws.Range("L49").Value = 17
$ws.Range("L54").Value = 61
$ws.Range("L55").Value = 31
$ws.Range("I63").Value = 246
$ws.Range("K63").Value = 79
$ws.Range("L63").Value = 16
$ws.Range("L64").Value = 25
$ws.Range("L65").Value = 60
$ws.Range("L67").Value = 101
$ws.Range("L71").Value = 8
$ws.Range("L73").Value = 20
$ws.Range("L75").Value = 13
$ws.Range("L76").Value = 39
$ws.Range("L78").Value = 42
$ws.Range("L79").Value = 84
$ws.Range("L81").Value = 4
$ws.Range("L83").Value = 59
$ws.Range("L84").Value = 26
$ws.Range("L85").Value = 144
$ws.Range("L86").Value = 20
$ws.Range("L89").Value = 36
$ws.Range("L93").Value = 13
$ws.Range("L95").Value = 43
$ws.Range("L96").Value = 25
$ws.Range("L99").Value = 46
$ws.Range("I101").Value = 26287
$ws.Range("K101").Value = 27526
$ws.Range("L101").Value = 2883

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 837
$ws.Range("L3").Value = 858
$ws.Range("I4").Value = 1823
$ws.Range("K4").Value = 1736
$ws.Range("L4").Value = 219
$ws.Range("L5").Value = 58
$ws.Range("L6").Value = 911
$ws.Range("I7").Value = 26287
$ws.Range("K7").Value = 27526
$ws.Range("L7").Value = 2883

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 34
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 95

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 50
$ws.Range("L3").Value = 56
$ws.Range("L6").Value = 59
$ws.Range("L7").Value = 181

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L2").Value = 3
$ws.Range("L3").Value = 8
$ws.Range("L7").Value = 18

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L2").Value = 8
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 94

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 76

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L3").Value = 2
$ws.Range("L7").Value = 8

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L7").Value = 8
$ws.Range("L4").Value = 1

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L2").Value = 2
$ws.Range("L7").Value = 5

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 46
$ws.Range("L6").Value = 38
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 29
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 125

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 30
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 95

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 24
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 49

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 14
$ws.Range("L4").Value = 2
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 25

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 19
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 101

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 8

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 13

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 4
$ws.Range("L4").Value = 10
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 13
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 42

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L5").Value = 4
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 84

$ws = $wb.Worksheets.Item("Sauganash,Forest Glen")
$ws.Range("L2").Value = 3
$ws.Range("L7").Value = 4
$ws.Range("L6").Value = 1

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 34
$ws.Range("L3").Value = 66
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 144

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 14
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L2").Value = 6
$ws.Range("L7").Value = 13

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L2").Value = 13
$ws.Range("L7").Value = 25

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 46
